$d = $word.ActiveDocument

# The author removed the word "first " from the sentence introducing the
# core modules list, turning:
#   "We first need to program its core modules which will be:"
# into:
#   "We need to program its core modules which will be:"
$d.Content.Find.Execute("We first need to program its core modules which will be:", $false, $false, $false, $false, $false, $true, 1, $false, "We need to program its core modules which will be:", 2)
